$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("message")

# Copy the formatting (styles) of the last existing data row (47) down into the
# five new rows (48-52) so the new rows look the same as the rest of the table.
$ws.Range("A47:C47").Copy()
$ws.Range("A48:C52").PasteSpecial(-4122)
$ws.Rows("48:52").RowHeight = 20

# Column A: running row number, same "=ROW()-2" formula used by the rest of the column.
# Assigning the whole block in one go lets the five cells form a single formula group.
$ws.Range("A48:A52").Formula = "=ROW()-2"

# Column B: the new item/ability messages being added.
$ws.Range("B48").Value = "<val1>の最大HPが<val2>上昇した"
$ws.Range("B49").Value = "<val1>の力が<val2>上昇した"
$ws.Range("B50").Value = "<val1>の体力が<val2>上昇した"
$ws.Range("B51").Value = "<val1>の素早さが<val2>上昇した"
$ws.Range("B52").Value = "<val1>の魔力が<val2>上昇した"

# Column C: color tag for these new messages, matching existing "green" rows.
$ws.Range("C48:C52").Value = "green"
